$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "1.003",
# "28.006.87") but must remain plain text, exactly as authored. Excel's
# COM layer auto-converts numeric-looking strings assigned via .Value
# into real numbers, so force the whole column to Text format before
# writing, and restore the original (default) cell style afterwards so
# no stray formatting is introduced.
$priceRange = $ws.Range("D2:D51")
$defaultStyle = $ws.Range("A2").Style
$priceRange.NumberFormat = "@"

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "28.006.87"
$ws.Range("E2").Value = "  -1.29%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.764.34"
$ws.Range("E3").Value = "  -3.43%  "

# Row 4 (TetherUSD)
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +1.19%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "321.16"
$ws.Range("E5").Value = "  -1.62%  "

# Row 6 (USDC)
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +1.22%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "0.4249"
$ws.Range("E7").Value = "  -4.75%  "

# Row 8 (Cardano)
$ws.Range("D8").Value = "0.3604"
$ws.Range("E8").Value = "  -4.99%  "

# Row 9 (OKB)
$ws.Range("D9").Value = "43.90"
$ws.Range("E9").Value = "  -2.56%  "

# Row 10 (Dogecoin)
$ws.Range("D10").Value = "0.07443"
$ws.Range("E10").Value = "  -4.46%  "

# Row 11 (Polygon)
$ws.Range("D11").Value = "1.101"
$ws.Range("E11").Value = "  -3.63%  "

# Row 12 (BinanceUSD)
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +1.27%  "

# Row 13 (Solana)
$ws.Range("D13").Value = "21.17"
$ws.Range("E13").Value = "  -4.96%  "

# Row 14 (Polkadot)
$ws.Range("D14").Value = "6.084"
$ws.Range("E14").Value = "  -3.46%  "

# Row 15 (Chainlink)
$ws.Range("D15").Value = "7.316"
$ws.Range("E15").Value = "  -2.81%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "1.786.94"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17 (Litecoin)
$ws.Range("D17").Value = "90.94"
$ws.Range("E17").Value = "  -1.64%  "

# Row 18 (ShibaInu)
$ws.Range("D18").Value = "0.00001057"
$ws.Range("E18").Value = "  -2.55%  "

# Row 19 (TRON)
$ws.Range("D19").Value = "0.06395"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20 (Dai)
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.97%  "

# Row 21 (Avalanche)
$ws.Range("D21").Value = "17.03"
$ws.Range("E21").Value = "  -3.33%  "

# Row 22 (Uniswap)
$ws.Range("D22").Value = "5.971"
$ws.Range("E22").Value = "  -6.19%  "

# Row 23 (WrappedBTC)
$ws.Range("D23").Value = "28.008.02"
$ws.Range("E23").Value = "  -1.29%  "

# Row 24 (Cosmos) - only E changes
$ws.Range("E24").Value = "  -3.66%  "

# Row 25 (Toncoin)
$ws.Range("D25").Value = "2.138"
$ws.Range("E25").Value = "  +3.67%  "

# Row 26 (Monero)
$ws.Range("D26").Value = "158.32"
$ws.Range("E26").Value = "  +3.05%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").Value = "20.17"
$ws.Range("E27").Value = "  -4.04%  "

# Row 28 (WrappedliquidstakedEther2.0)
$ws.Range("D28").Value = "1.985.83"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29 (LidoDAOToken)
$ws.Range("D29").Value = "2.129"
$ws.Range("E29").Value = "  -10.63%  "

# Row 30 (BitcoinCash)
$ws.Range("D30").Value = "124.93"
$ws.Range("E30").Value = "  -3.64%  "

# Row 31 (ImmutableX)
$ws.Range("D31").Value = "1.166"
$ws.Range("E31").Value = "  -4.42%  "

# Row 32 (Filecoin)
$ws.Range("D32").Value = "5.642"
$ws.Range("E32").Value = "  -3.79%  "

# Row 33 (Stellar)
$ws.Range("D33").Value = "0.08870"
$ws.Range("E33").Value = "  -4.05%  "

# Row 34 (HuobiToken)
$ws.Range("D34").Value = "3.537"
$ws.Range("E34").Value = "  -3.02%  "

# Row 35 (Aptos)
$ws.Range("D35").Value = "12.54"
$ws.Range("E35").Value = "  -2.42%  "

# Row 36 (VeChain)
$ws.Range("D36").Value = "0.02316"
$ws.Range("E36").Value = "  -1.75%  "

# Row 37 (Algorand) - only E changes
$ws.Range("E37").Value = "  -4.27%  "

# Row 38 / 39: InternetComputer(DFINITY) and Hedera swap places
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06051"
$ws.Range("E38").Value = "  -2.79%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.022"
$ws.Range("E39").Value = "  -3.31%  "

# Row 40 (TheSandbox)
$ws.Range("D40").Value = "0.6371"
$ws.Range("E40").Value = "  -3.86%  "

# Row 41 (TrustWalletToken)
$ws.Range("D41").Value = "1.184"
$ws.Range("E41").Value = "  -0.81%  "

# Row 42 (Frax)
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +1.15%  "

# Row 43 / 44: FraxShare and WEMIXTOKEN swap places
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "1.398"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "7.836"
$ws.Range("E44").Value = "  -3.42%  "

# Row 45 (EnergySwap)
$ws.Range("D45").Value = "13.53"
$ws.Range("E45").Value = "  -2.90%  "

# Row 46 (Decentraland)
$ws.Range("D46").Value = "0.5918"
$ws.Range("E46").Value = "  -3.51%  "

# Row 47 (PancakeSwap)
$ws.Range("D47").Value = "3.695"
$ws.Range("E47").Value = "  -1.26%  "

# Row 48 (NEARProtocol)
$ws.Range("D48").Value = "2.005"
$ws.Range("E48").Value = "  -1.57%  "

# Row 49 (Quant)
$ws.Range("D49").Value = "123.08"
$ws.Range("E49").Value = "  -3.32%  "

# Row 50 (EOS)
$ws.Range("D50").Value = "1.194"
$ws.Range("E50").Value = "  +4.10%  "

# Row 51 (Cronos)
$ws.Range("D51").Value = "0.06869"
$ws.Range("E51").Value = "  -1.97%  "

# Restore original (unstyled) formatting for the Price column now that
# all text values have been written.
$priceRange.Style = $defaultStyle
